# Generate Report for Handback
# Refresh the handback timestamps for the second tracked file
# (82d82f6b-1ec0-4ac5-a4f0-77bc32a55d74) across the Overview, zh-cn and
# de-de report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-28 16:48:57"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-28 16:48:53"
$zhcn.Range("K3").Value = "2016-08-28 16:49:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-28 16:48:57"
$dede.Range("K3").Value = "2016-08-28 16:49:19"
